$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..21 down to 3..22)
$ws.Rows(2).Insert()
$ws.Range("A2:F2").ClearFormats()

# Populate the new row 2 with the 에이피알 IPO record
$ws.Cells.Item(2, 1).Value = "에이피알"
$ws.Cells.Item(2, 2).Value = "2024.01.22~01.26"
$ws.Cells.Item(2, 3).Value = "147,000~200,000"
$ws.Cells.Item(2, 4).Value = "-"
$ws.Cells.Item(2, 5).Value = 55713
$ws.Cells.Item(2, 6).Value = "신한투자증권,하나증권"

# Drop the oldest record that was pushed past row 21 (에이텀, 2023.11.09~11.15 entry)
$ws.Rows(22).Delete()
